# Applies the "Updated cryptos list" data refresh to Sheet1.
# Plain text/string updates (coin names, links, percentage-change text)
# and numeric-looking price updates that must remain TEXT cells (as in
# the source data, which uses dotted thousands separators like "65.609.06").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates whose new text is unambiguously non-numeric: plain assignment is safe ---
$plainUpdates = @(
    @{Cell='D2'; Value='65.609.06'},
    @{Cell='E2'; Value='  -3.53%  '},
    @{Cell='D3'; Value='3.477.71'},
    @{Cell='E3'; Value='  -0.24%  '},
    @{Cell='E5'; Value='  -2.25%  '},
    @{Cell='E6'; Value='  -5.08%  '},
    @{Cell='E7'; Value='  +0.01%  '},
    @{Cell='E8'; Value='  -4.10%  '},
    @{Cell='D9'; Value='3.474.83'},
    @{Cell='E9'; Value='  -0.29%  '},
    @{Cell='E10'; Value='  -7.64%  '},
    @{Cell='E12'; Value='  -5.12%  '},
    @{Cell='D13'; Value='4.077.54'},
    @{Cell='E13'; Value='  -0.25%  '},
    @{Cell='E14'; Value='  +0.07%  '},
    @{Cell='E15'; Value='  -6.80%  '},
    @{Cell='D16'; Value='65.816.13'},
    @{Cell='E16'; Value='  -3.28%  '},
    @{Cell='E17'; Value='  -4.00%  '},
    @{Cell='D18'; Value='3.476.54'},
    @{Cell='E18'; Value='  -0.46%  '},
    @{Cell='E19'; Value='  -4.54%  '},
    @{Cell='E20'; Value='  -1.54%  '},
    @{Cell='E21'; Value='  -7.51%  '},
    @{Cell='E22'; Value='  -2.46%  '},
    @{Cell='E23'; Value='  +0.05%  '},
    @{Cell='E24'; Value='  +0.12%  '},
    @{Cell='E25'; Value='  -1.04%  '},
    @{Cell='E26'; Value='  -0.97%  '},
    @{Cell='E27'; Value='  -7.25%  '},
    @{Cell='E28'; Value='  +0.26%  '},
    @{Cell='E29'; Value='  +0.25%  '},
    @{Cell='E30'; Value='  +1.90%  '},
    @{Cell='E31'; Value='  -6.22%  '},
    @{Cell='E32'; Value='  -4.13%  '},
    @{Cell='E33'; Value='  -0.04%  '},
    @{Cell='E34'; Value='  -3.61%  '},
    @{Cell='E35'; Value='  -8.60%  '},
    @{Cell='E36'; Value='  -2.14%  '},
    @{Cell='B37'; Value='Monero'},
    @{Cell='C37'; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'},
    @{Cell='E37'; Value='  -1.25%  '},
    @{Cell='B38'; Value='EnergySwap'},
    @{Cell='C38'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'},
    @{Cell='E38'; Value='  +11.35%  '},
    @{Cell='E39'; Value='  -0.65%  '},
    @{Cell='E40'; Value='  -5.57%  '},
    @{Cell='D41'; Value='2.796.29'},
    @{Cell='E41'; Value='  +1.41%  '},
    @{Cell='E42'; Value='  -9.75%  '},
    @{Cell='E43'; Value='  -5.21%  '},
    @{Cell='E44'; Value='  -6.98%  '},
    @{Cell='E45'; Value='  -4.57%  '},
    @{Cell='E46'; Value='  -3.94%  '},
    @{Cell='E47'; Value='  -8.81%  '},
    @{Cell='E48'; Value='  -3.88%  '},
    @{Cell='E49'; Value='  -6.83%  '},
    @{Cell='E50'; Value='  -3.10%  '},
    @{Cell='B51'; Value='Cosmos'},
    @{Cell='C51'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'},
    @{Cell='E51'; Value='  -2.62%  '}
)

foreach ($u in $plainUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Updates whose new text LOOKS like a plain number (e.g. "30.00", "1.00",
#     "0.0000121"): Excel auto-converts such literals to numeric cells, which
#     would silently drop the original text formatting (trailing zeros, etc.)
#     and change the stored cell type from string to number. Force these to
#     stay text by flipping to a text number-format for the write, then restore
#     the cell's original style so no formatting footprint is left behind. ---
$textUpdates = @(
    @{Cell='D5'; Value='581.77'},
    @{Cell='D6'; Value='172.69'},
    @{Cell='D15'; Value='30.00'},
    @{Cell='D19'; Value='5.94'},
    @{Cell='D21'; Value='366.12'},
    @{Cell='D22'; Value='7.77'},
    @{Cell='D23'; Value='1.00'},
    @{Cell='D24'; Value='72.38'},
    @{Cell='D25'; Value='0.534'},
    @{Cell='D26'; Value='0.0000121'},
    @{Cell='D30'; Value='24.09'},
    @{Cell='D34'; Value='7.11'},
    @{Cell='D37'; Value='159.97'},
    @{Cell='D38'; Value='29.28'},
    @{Cell='D40'; Value='1.78'},
    @{Cell='D42'; Value='2.56'},
    @{Cell='D44'; Value='6.32'},
    @{Cell='D46'; Value='39.87'},
    @{Cell='D47'; Value='24.08'},
    @{Cell='D49'; Value='307.83'},
    @{Cell='D51'; Value='6.21'}
)

foreach ($u in $textUpdates) {
    $range = $ws.Range($u.Cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $u.Value
    $range.Style = $origStyle
}

Write-Output "Applied $($plainUpdates.Count + $textUpdates.Count) cell updates"
